$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "3pm - "
$ws.Range("D28").Value = "Need to: make layout for forum, thread"
$ws.Range("D29").Value = "Need to decide about abilities and presentation for admin"
$ws.Range("D30").Value = "Need to revise intro and probably user stories to reflect evolution here."

$ws.Range("D35").Select()
